$wb = $excel.ActiveWorkbook

# --- CSV sheet content updates ---
$csv = $wb.Worksheets.Item("CSV")

# Fix missing opening quote in the "All - IR" query (B1)
$csv.Range("B1").Value = 'project = CI AND issuetype in (standardIssueTypes(), "Expense Delivery") AND "Epic Link" is EMPTY AND "Case Number/s" is not EMPTY AND cf[14125] in ("Incident Request (IR)") AND resolved is not EMPTY AND resolutiondate >= 2022-12-19'

# Replace the "Individual Tasks" query (B5) with the new JQL text
$csv.Range("B5").Value = 'project = CI AND  type = "Individual Task" AND resolveDate >= 2023-07-03 AND resolveDate  <= 2023-08-17'

# New stray cell content
$csv.Range("B12").Value = "c"

# --- Insert a new blank "Sheet1" between InProgress and SLA ---
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "Sheet1"
$sla = $wb.Worksheets.Item("SLA")
$newSheet.Move($sla)

# Update selection on the CSV sheet (after sheet insert, so indices are stable)
[void]$csv.Range("B5").Select()
